# Fix contact information missing from short resumes:
# Add a centered contact-info paragraph right after the "Dheeraj Chand"
# name/title paragraph (mirrors what the long-resume template already has).
#
# We do this via Find/Replace rather than Paragraphs.InsertParagraphAfter/
# Before so the new paragraph/run don't inherit the name run's character
# formatting (bold, 28pt) or a neighboring heading's paragraph style -
# the replacement text picks up the plain paragraph/run defaults, matching
# how the new line only carries the <w:jc w:val="center"/> alignment.

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Dheeraj Chand",
    $false,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "Dheeraj Chand^p202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX",
    2
)
